# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" on Overview sheet, and the
# corresponding "Priority" / "Latest Handoff Datetime" columns on the
# zh-cn and de-de detail sheets for the six files that were just handed
# off (rows 4,5,6,7,8,10 -- row 9 was already handed off earlier and is
# left untouched).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(4, 5, 6, 7, 8, 10)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-10-18 05:29:02"

    # zh-cn detail sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-10-18 05:28:46"

    # de-de detail sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-10-18 05:29:02"
}
